$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C, matching style of existing header cells (A1/B1):
# bold font, thin box border, horizontally centered / top-aligned.
$ws.Range("C1").Value = "min_units"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160

# Row 2 -> 3 units
$ws.Cells.Item(2, 3).Value = 3

# Rows 3-7 -> 4.5 units
for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = 4.5
}

# Rows 8-289 -> 9 units
for ($r = 8; $r -le 289; $r++) {
    $ws.Cells.Item($r, 3).Value = 9
}
